$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row at 57 (pushes the existing rows 57-61 down to 58-62) and
# populate it with the missing "X5" 32.768kHz crystal (ABS07) that this
# commit adds to the BOM.
# ---------------------------------------------------------------------------
$ws.Rows.Item(57).Insert()

$ws.Range("A57").Value = 1
$ws.Range("B57").Value = "X5"
$ws.Range("C57").Value = "32.768kHz"
$ws.Range("D57").Value = "N"
$ws.Range("F57").Value = "ABS07-32.768KHZ-7-T"
$ws.Range("G57").Value = "ABS07"
$ws.Range("H57").Value = "32.768kHz Crystal"
$ws.Range("I57").Value = "CRYSTAL 32.7680KHZ 7PF SMD"
$ws.Range("J57").Value = "535-9543-1-ND"
$ws.Range("K57").Value = "Abracon LLC"
$ws.Range("L57").Value = "ABS07-32.768KHZ-7-T"

# ---------------------------------------------------------------------------
# Conditional formatting: the two ranges that only covered rows 59:61 before
# the insert need to track the same rows, now 60:62.
# ---------------------------------------------------------------------------
$fcD = $ws.Range("D59:D61").FormatConditions.Item(1)
$fcD.ModifyAppliesToRange($ws.Range("D60:D62"))

$fcE = $ws.Range("E59:E61").FormatConditions.Item(1)
$fcE.ModifyAppliesToRange($ws.Range("E60:E62"))

# The first block (rows 53:58 unioned with D2:D10/D34:D50/D12:D29 etc.) needs
# to keep covering the crystal rows, which now run through row 59. Patch the
# newly-uncovered row with matching rules (same colours as the existing
# dxfs, OLE BGR-packed) rather than rewriting the whole multi-area rule.
$fcY59 = $ws.Range("D59").FormatConditions.Add(1, 3, '"Y"')
$fcY59.Font.Color = 26112
$fcY59.Interior.Color = 13434828
$fcY59.SetLastPriority()

$fcN59 = $ws.Range("D59").FormatConditions.Add(1, 3, '"N"')
$fcN59.Font.Color = 204
$fcN59.Interior.Color = 13421823
$fcN59.SetLastPriority()

$fcDnp59 = $ws.Range("E59").FormatConditions.Add(1, 3, '"DNP"')
$fcDnp59.Font.Color = 16777215
$fcDnp59.Font.Bold = $true
$fcDnp59.Interior.Color = 204
$fcDnp59.SetLastPriority()

# ---------------------------------------------------------------------------
# View state: zoom in to 90% and leave the selection on L58 (matches what a
# reviewer would land on after checking the new row's MPN cell).
# ---------------------------------------------------------------------------
$ws.Range("L58").Select() | Out-Null
$excel.ActiveWindow.Zoom = 90
